# Factor out resolution_z_unit: add "mm" as an allowed unit value.
#
# This mirrors the upstream change which:
#   1. Adds a new shared string "mm".
#   2. Updates the "resolution_z_unit list" sheet (sheet7) so it now
#      contains three values: mm, um, nm (in that order).
#   3. Updates the "Export as TSV" sheet's data validation on the
#      resolution_z_unit column (V) so that it references the 3-row
#      range and reports the new set of allowed values in its error
#      message.

$wb = $excel.ActiveWorkbook

# 1 & 2. Update the "resolution_z_unit list" sheet.
#    Before: A1 = nm, A2 = um
#    After:  A1 = mm, A2 = um, A3 = nm
$wsZ = $wb.Worksheets.Item("resolution_z_unit list")
$wsZ.Range("A3").Value = $wsZ.Range("A1").Value2   # A3 = nm (moved down)
$wsZ.Range("A1").Value = "mm"                      # A1 = mm (new)
# A2 already holds "um" and stays untouched.

# 3. Update the data validation on the resolution_z_unit column (V) of
#    the main "Export as TSV" sheet so it points at the expanded list
#    and reports the new allowed values.
$wsMain = $wb.Worksheets.Item("Export as TSV")
$dv = $wsMain.Range("V2:V1048576").Validation
$dv.Delete()
$dv.Add(3, 1, 1, "='resolution_z_unit list'!`$A`$1:`$A`$3")
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: mm / um / nm."
$dv.ShowInput = $true
$dv.ShowError = $true
